$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "30.415.31"
Set-TextValue $ws.Range("E2") "  +2.50%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.109.44"
Set-TextValue $ws.Range("E3") "  +0.64%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.006"
Set-TextValue $ws.Range("E4") "  -0.23%  "

# Row 5
Set-TextValue $ws.Range("D5") "345.24"
Set-TextValue $ws.Range("E5") "  +0.75%  "

# Row 6
Set-TextValue $ws.Range("E6") "  -0.17%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.5234"
Set-TextValue $ws.Range("E7") "  +2.09%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.4453"
Set-TextValue $ws.Range("E8") "  +1.21%  "

# Row 9
Set-TextValue $ws.Range("D9") "54.82"
Set-TextValue $ws.Range("E9") "  +3.08%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.09388"
Set-TextValue $ws.Range("E10") "  +2.56%  "

# Row 11
Set-TextValue $ws.Range("E11") "  +0.46%  "

# Row 12
Set-TextValue $ws.Range("D12") "25.00"
Set-TextValue $ws.Range("E12") "  +0.91%  "

# Row 13
Set-TextValue $ws.Range("D13") "8.702"
Set-TextValue $ws.Range("E13") "  +6.31%  "

# Row 14
Set-TextValue $ws.Range("D14") "6.957"
Set-TextValue $ws.Range("E14") "  +3.27%  "

# Row 15
Set-TextValue $ws.Range("D15") "2.019.50"
Set-TextValue $ws.Range("E15") "  -4.11%  "

# Row 16
Set-TextValue $ws.Range("D16") "102.12"
Set-TextValue $ws.Range("E16") "  +2.51%  "

# Row 17
Set-TextValue $ws.Range("D17") "0.00001163"
Set-TextValue $ws.Range("E17") "  +1.24%  "

# Row 18
Set-TextValue $ws.Range("E18") "  -0.07%  "

# Row 19
Set-TextValue $ws.Range("D19") "21.27"
Set-TextValue $ws.Range("E19") "  +0.97%  "

# Row 20
Set-TextValue $ws.Range("E20") "  +1.23%  "

# Row 21
Set-TextValue $ws.Range("D21") "6.365"
Set-TextValue $ws.Range("E21") "  +3.20%  "

# Row 22
Set-TextValue $ws.Range("D22") "1.006"

# Row 23
Set-TextValue $ws.Range("D23") "30.444.32"
Set-TextValue $ws.Range("E23") "  +2.39%  "

# Row 24
Set-TextValue $ws.Range("D24") "12.67"
Set-TextValue $ws.Range("E24") "  +0.72%  "

# Row 25
Set-TextValue $ws.Range("D25") "2.298"
Set-TextValue $ws.Range("E25") "  -0.54%  "

# Row 26
Set-TextValue $ws.Range("E26") "  +1.16%  "

# Row 27
Set-TextValue $ws.Range("D27") "163.02"
Set-TextValue $ws.Range("E27") "  +0.22%  "

# Row 28
Set-TextValue $ws.Range("D28") "2.531"
Set-TextValue $ws.Range("E28") "  +0.75%  "

# Row 29
Set-TextValue $ws.Range("D29") "134.06"
Set-TextValue $ws.Range("E29") "  +1.21%  "

# Row 30
Set-TextValue $ws.Range("D30") "1.153"
Set-TextValue $ws.Range("E30") "  +2.37%  "

# Row 31
Set-TextValue $ws.Range("D31") "1.754"
Set-TextValue $ws.Range("E31") "  +7.38%  "

# Row 32
Set-TextValue $ws.Range("E32") "  +1.06%  "

# Row 33
Set-TextValue $ws.Range("D33") "6.837"
Set-TextValue $ws.Range("E33") "  +13.64%  "

# Row 34
Set-TextValue $ws.Range("D34") "6.275"
Set-TextValue $ws.Range("E34") "  +2.04%  "

# Row 35
Set-TextValue $ws.Range("D35") "3.922"
Set-TextValue $ws.Range("E35") "  -0.96%  "

# Row 36
Set-TextValue $ws.Range("D36") "10.35"
Set-TextValue $ws.Range("E36") "  +1.21%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.02629"
Set-TextValue $ws.Range("E37") "  +2.35%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.06798"
Set-TextValue $ws.Range("E38") "  +1.79%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.7062"
Set-TextValue $ws.Range("E39") "  +3.33%  "

# Row 40
Set-TextValue $ws.Range("E40") "  +5.47%  "

# Row 41
Set-TextValue $ws.Range("D41") "12.62"
Set-TextValue $ws.Range("E41") "  +2.05%  "

# Row 42
Set-TextValue $ws.Range("D42") "0.2230"
Set-TextValue $ws.Range("E42") "  +0.25%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.6873"
Set-TextValue $ws.Range("E43") "  +3.09%  "

# Row 44
Set-TextValue $ws.Range("D44") "14.54"
Set-TextValue $ws.Range("E44") "  +2.49%  "

# Row 45
Set-TextValue $ws.Range("D45") "2.362"
Set-TextValue $ws.Range("E45") "  +3.07%  "

# Row 46
Set-TextValue $ws.Range("E46") "  -0.04%  "

# Row 47
Set-TextValue $ws.Range("D47") "1.394"
Set-TextValue $ws.Range("E47") "  +20.11%  "

# Row 48
Set-TextValue $ws.Range("D48") "3.650"
Set-TextValue $ws.Range("E48") "  +1.27%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.00000000352"
Set-TextValue $ws.Range("E49") "  +5.56%  "

# Row 50
Set-TextValue $ws.Range("D50") "1.213"
Set-TextValue $ws.Range("E50") "  +9.87%  "

# Row 51
Set-TextValue $ws.Range("D51") "1.222"
Set-TextValue $ws.Range("E51") "  +0.33%  "
